$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column V: "Parcial Rodada 1" - partial score for round 1
# Copy the header style (bold, centered, bordered) from U1 into V1,
# matching the look of the other header cells in row 1.
$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("V1").Value = "Parcial Rodada 1"

# Populate "Parcial Rodada 1" partial scores for each team (V2:V33)
$ws.Range("V2").Value = 59.56
$ws.Range("V3").Value = 45.46
$ws.Range("V4").Value = 64.95999999999999
$ws.Range("V5").Value = 61.8
$ws.Range("V6").Value = 50.76
$ws.Range("V7").Value = 56.65
$ws.Range("V8").Value = 59.86
$ws.Range("V9").Value = 39.66
$ws.Range("V10").Value = 60.66
$ws.Range("V11").Value = 23.26
$ws.Range("V12").Value = 53.06
$ws.Range("V13").Value = 59.76
$ws.Range("V14").Value = 34.76
$ws.Range("V15").Value = 69.56
$ws.Range("V16").Value = 41.6
$ws.Range("V17").Value = 73.66
$ws.Range("V18").Value = 51.05
$ws.Range("V19").Value = 49.16
$ws.Range("V20").Value = 50.6
$ws.Range("V21").Value = 44.46
$ws.Range("V22").Value = 69.76000000000001
$ws.Range("V23").Value = 58.4
$ws.Range("V24").Value = 43.1
$ws.Range("V25").Value = 30.6
$ws.Range("V26").Value = 44.06
$ws.Range("V27").Value = 38.5
$ws.Range("V28").Value = 57.76
$ws.Range("V29").Value = 34.36
$ws.Range("V30").Value = 54.36
$ws.Range("V31").Value = 69.26000000000001
$ws.Range("V32").Value = 60.36
$ws.Range("V33").Value = 39.66
